# ============================================================
# Commit: Implementa suporte para múltiplas abas no processamento
# de planilhas Excel - adiciona aba "Ações", renomeia "Sheet1" -> "Consolidado"
# ============================================================

$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> Consolidado, insert new sheet "Ações" right after it ---
$consolidado = $wb.Worksheets.Item(1)
$consolidado.Name = "Consolidado"

$acoes = $wb.Worksheets.Add($null, $consolidado)
$acoes.Name = "Ações"

# --- Style #1 (wrapText) must be allocated first: applied to header C1 ---
$acoes.Range("C1").WrapText = $true

# --- Header row: original 9 columns (A, D..K) written first ---
# (this reproduces the shared-string table order: indices 56..64)
$acoes.Range("A1").Value = "Ticker"
$acoes.Range("D1").Value = "Renda Esperada"
$acoes.Range("E1").Value = "Capital Atual"
$acoes.Range("F1").Value = "Dividend Yield Esperado"
$acoes.Range("G1").Value = "Dividend Yield Pago"
$acoes.Range("H1").Value = "Proporção Hoje"
$acoes.Range("I1").Value = "Meta 28k"
$acoes.Range("J1").Value = "Meta +1.a."
$acoes.Range("K1").Value = "Meta qtd. 2033"

# --- Style #2 (0.00%) then Style #3 (#,##0): allocated before the number cells below ---
$acoes.Range("F2:F11").NumberFormat = "0.00%"
$acoes.Range("H2:H11").NumberFormat = "0.00%"
$acoes.Range("J2:J3").NumberFormat = "#,##0"
$acoes.Range("K2:K10").NumberFormat = "#,##0"

# --- Data rows 2..11 ---
# Row 2
$acoes.Range("A2").Value = "BBAS3"
$acoes.Range("B2").Value = 5037
$acoes.Range("C2").Value = "R`$1.10"
$acoes.Range("D2").Value = "R`$5,558.31"
$acoes.Range("E2").Value = "R`$104,014.05"
$acoes.Range("F2").Value = 0.0534
$acoes.Range("G2").Value = "R`$5,302.45"
$acoes.Range("H2").Value = 0.3362
$acoes.Range("I2").Value = -142
$acoes.Range("J2").Value = -2519
$acoes.Range("K2").Value = -45417

# Row 3
$acoes.Range("A3").Value = "BBSE3"
$acoes.Range("B3").Value = 3220
$acoes.Range("C3").Value = "R`$3.81"
$acoes.Range("D3").Value = "R`$12,261.52"
$acoes.Range("E3").Value = "R`$105,583.80"
$acoes.Range("F3").Value = 0.1161
$acoes.Range("G3").Value = "R`$7,316.81"
$acoes.Range("H3").Value = 0.3413
$acoes.Range("I3").Value = -91
$acoes.Range("J3").Value = -1610
$acoes.Range("K3").Value = -29034

# Row 4
$acoes.Range("A4").Value = "BRSR6"
$acoes.Range("B4").Value = 1314
$acoes.Range("C4").Value = "R`$1.10"
$acoes.Range("D4").Value = "R`$1,445.81"
$acoes.Range("E4").Value = "R`$15,531.48"
$acoes.Range("F4").Value = 0.0931
$acoes.Range("G4").Value = "R`$692.74"
$acoes.Range("H4").Value = 0.0502
$acoes.Range("I4").Value = -37
$acoes.Range("J4").Value = -657
$acoes.Range("K4").Value = -11848

# Row 5
$acoes.Range("A5").Value = "CEBR6"
$acoes.Range("B5").Value = 645
$acoes.Range("C5").Value = "R`$2.29"
$acoes.Range("D5").Value = "R`$1,476.44"
$acoes.Range("E5").Value = "R`$15,480.00"
$acoes.Range("F5").Value = 0.0954
$acoes.Range("G5").Value = "R`$815.67"
$acoes.Range("H5").Value = 0.05
$acoes.Range("I5").Value = -18
$acoes.Range("J5").Value = -323
$acoes.Range("K5").Value = -5816

# Row 6
$acoes.Range("A6").Value = "CXSE3"
$acoes.Range("B6").Value = 1410
$acoes.Range("C6").Value = "R`$1.32"
$acoes.Range("D6").Value = "R`$1,861.20"
$acoes.Range("E6").Value = "R`$19,204.20"
$acoes.Range("F6").Value = 0.0969
$acoes.Range("G6").Value = "R`$1,325.40"
$acoes.Range("H6").Value = 0.0621
$acoes.Range("I6").Value = -40
$acoes.Range("J6").Value = -705
$acoes.Range("K6").Value = -12714

# Row 7
$acoes.Range("A7").Value = "LEVE3"
$acoes.Range("B7").Value = 634
$acoes.Range("C7").Value = "R`$2.74"
$acoes.Range("D7").Value = "R`$1,735.40"
$acoes.Range("E7").Value = "R`$17,396.96"
$acoes.Range("F7").Value = 0.0998
$acoes.Range("G7").Value = "R`$1,329.75"
$acoes.Range("H7").Value = 0.0562
$acoes.Range("I7").Value = -18
$acoes.Range("J7").Value = -317
$acoes.Range("K7").Value = -5717

# Row 8
$acoes.Range("A8").Value = "PETR4"
$acoes.Range("B8").Value = 289
$acoes.Range("C8").Value = "R`$2.97"
$acoes.Range("D8").Value = "R`$857.67"
$acoes.Range("E8").Value = "R`$8,719.13"
$acoes.Range("F8").Value = 0.0984
$acoes.Range("G8").Value = "R`$588.92"
$acoes.Range("H8").Value = 0.0282
$acoes.Range("I8").Value = -8
$acoes.Range("J8").Value = -145
$acoes.Range("K8").Value = -2606

# Row 9
$acoes.Range("A9").Value = "RANI3"
$acoes.Range("B9").Value = 1460
$acoes.Range("C9").Value = "R`$0.76"
$acoes.Range("D9").Value = "R`$1,109.61"
$acoes.Range("E9").Value = "R`$10,818.60"
$acoes.Range("F9").Value = 0.1026
$acoes.Range("G9").Value = "R`$845.49"
$acoes.Range("H9").Value = 0.035
$acoes.Range("I9").Value = -41
$acoes.Range("J9").Value = -730
$acoes.Range("K9").Value = -13164

# Row 10
$acoes.Range("A10").Value = "ISAE4"
$acoes.Range("B10").Value = 496
$acoes.Range("C10").Value = "R`$1.59"
$acoes.Range("D10").Value = "R`$790.43"
$acoes.Range("E10").Value = "R`$10,976.48"
$acoes.Range("F10").Value = 0.072
$acoes.Range("G10").Value = "R`$1,170.91"
$acoes.Range("H10").Value = 0.0355
$acoes.Range("I10").Value = -14
$acoes.Range("J10").Value = -248
$acoes.Range("K10").Value = -4472

# Row 11
$acoes.Range("A11").Value = "CGAS5"
$acoes.Range("B11").Value = 13
$acoes.Range("C11").Value = "R`$10.56"
$acoes.Range("D11").Value = "R`$137.34"
$acoes.Range("E11").Value = "R`$1,664.00"
$acoes.Range("F11").Value = 0.0825
$acoes.Range("G11").Value = "R`$9.95"
$acoes.Range("H11").Value = 0.0054
$acoes.Range("I11").Value = 0
$acoes.Range("J11").Value = -7
$acoes.Range("K11").Value = -117

# --- Finally, the two header cells added last by the original author ---
# (this is why "Div. Esperado 2025" / "Qtd" are the LAST two shared strings)
$acoes.Range("C1").Value = "Div. Esperado 2025"
$acoes.Range("B1").Value = "Qtd"

# --- Column widths (approximate bestFit) ---
$acoes.Columns("A").ColumnWidth = 6.7109375
$acoes.Columns("B").ColumnWidth = 6.28515625
$acoes.Columns("C").ColumnWidth = 18.140625
$acoes.Columns("D").ColumnWidth = 15.28515625
$acoes.Columns("E").ColumnWidth = 13.140625
$acoes.Columns("F").ColumnWidth = 22.85546875
$acoes.Columns("G").ColumnWidth = 18.85546875
$acoes.Columns("H").ColumnWidth = 14.85546875
$acoes.Columns("J").ColumnWidth = 10.140625
$acoes.Columns("K").ColumnWidth = 14.28515625

# --- View / selection state ---
$acoes.Range("G1").Select()

# --- Restore Consolidado sheet selection (B32) ---
$consolidado.Range("B32").Select()
$acoes.Select()

